$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Vanern (row 11) values ---
# I11 keeps its existing style (s="6"); just set the value.
$ws.Range("I11").Value = 4680000

# J11/K11 in the target no longer carry their old style index (s="3")
# — reset to Normal before assigning the new values.
$ws.Range("J11").Style = "Normal"
$ws.Range("K11").Style = "Normal"
$ws.Range("J11").Value = 0.037
$ws.Range("K11").Value = 0.6147

# --- Update Vanern source/comment text (B24) ---
# Replace the old placeholder comment with the real citation. Since this
# was the only cell referencing that shared string, the stale string is
# dropped and the new one appended on save.
$ws.Range("B24").Value = "Kvarnas, H (2001) Morphometry and hydrology of the four large lakes of Sweden. Ambio 30(8): 467-474; other sources see CommentsQuirks file for links."

# --- Update active selection to match the saved view state ---
$ws.Range("B25").Select()
